# Apply the manuscript-figure update:
#  - Column A header "paper" -> "study"
#  - Study name "McCoy et al"    -> "McCoy et al 2021"
#  - Study name "Kaitany et al"  -> "Kaitany et al 2001"
#  - Selection moved to E5 (cosmetic, matches the saved view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row
$ws.Range("A1").Value = "study"

# Update the two study-name columns (column A, rows 2-29)
for ($r = 2; $r -le 29; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "McCoy et al") {
        $cell.Value = "McCoy et al 2021"
    } elseif ($cell.Value2 -eq "Kaitany et al") {
        $cell.Value = "Kaitany et al 2001"
    }
}

# Restore the active selection recorded in the saved view
$ws.Range("E5").Select()
